# Rename "Gorillas" to "Golliras" on Sheet1, and update the selected cell
# to B4 (matching the author's final selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A3 currently holds "Gorillas" -> rename to "Golliras"
$ws.Range("A3").Value = "Golliras"

# Update the active selection to B4, as recorded in the saved view state
$ws.Activate()
$ws.Range("B4").Select()
